$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44358
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 14000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 14500
$ws.Cells.Item(2, 16).Value = 1450

# Row 3
$ws.Cells.Item(3, 4).Value = 44943
$ws.Cells.Item(3, 10).Value = 400
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 17000
$ws.Cells.Item(3, 13).Value = 16500
$ws.Cells.Item(3, 16).Value = 1650

# Row 4
$ws.Cells.Item(4, 4).Value = 44860
$ws.Cells.Item(4, 10).Value = 400
$ws.Cells.Item(4, 11).Value = 14000
$ws.Cells.Item(4, 12).Value = 15000
$ws.Cells.Item(4, 13).Value = 14500
$ws.Cells.Item(4, 16).Value = 1450

# Row 5
$ws.Cells.Item(5, 4).Value = 44890
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 17000
$ws.Cells.Item(5, 13).Value = 16500
$ws.Cells.Item(5, 16).Value = 1650

# Row 6
$ws.Cells.Item(6, 4).Value = 45001
$ws.Cells.Item(6, 10).Value = 400
$ws.Cells.Item(6, 11).Value = 18000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 19000
$ws.Cells.Item(6, 16).Value = 1900

# Row 7
$ws.Cells.Item(7, 4).Value = 44160
$ws.Cells.Item(7, 10).Value = 360
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 11000
$ws.Cells.Item(7, 13).Value = 10500
$ws.Cells.Item(7, 16).Value = 1050

# Row 8
$ws.Cells.Item(8, 4).Value = 44377
$ws.Cells.Item(8, 10).Value = 650
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 14538
$ws.Cells.Item(8, 16).Value = 1454

# Row 9
$ws.Cells.Item(9, 4).Value = 44330
$ws.Cells.Item(9, 10).Value = 300
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 13500
$ws.Cells.Item(9, 16).Value = 1350

# Row 10
$ws.Cells.Item(10, 4).Value = 44942
$ws.Cells.Item(10, 10).Value = 1000
$ws.Cells.Item(10, 11).Value = 14000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 14500
$ws.Cells.Item(10, 16).Value = 1450

# Row 11
$ws.Cells.Item(11, 4).Value = 44972
$ws.Cells.Item(11, 10).Value = 550
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 16000
$ws.Cells.Item(11, 13).Value = 15636
$ws.Cells.Item(11, 16).Value = 1564

# Row 12
$ws.Cells.Item(12, 4).Value = 44204
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 16).Value = 1050

# Row 13
$ws.Cells.Item(13, 4).Value = 44959
$ws.Cells.Item(13, 10).Value = 400
$ws.Cells.Item(13, 11).Value = 21000
$ws.Cells.Item(13, 12).Value = 22000
$ws.Cells.Item(13, 13).Value = 21500
$ws.Cells.Item(13, 16).Value = 2150

# Row 14
$ws.Cells.Item(14, 4).Value = 44914
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = 14500
$ws.Cells.Item(14, 16).Value = 1450

# Row 15
$ws.Cells.Item(15, 4).Value = 44406
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 14500
$ws.Cells.Item(15, 16).Value = 1450

# Row 16
$ws.Cells.Item(16, 4).Value = 44547
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 19000
$ws.Cells.Item(16, 12).Value = 20000
$ws.Cells.Item(16, 13).Value = 19500
$ws.Cells.Item(16, 16).Value = 1950

# Row 17
$ws.Cells.Item(17, 4).Value = 44893
$ws.Cells.Item(17, 10).Value = 1400
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 16000
$ws.Cells.Item(17, 13).Value = 15571
$ws.Cells.Item(17, 16).Value = 1557

# Row 18
$ws.Cells.Item(18, 4).Value = 44291
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 13500
$ws.Cells.Item(18, 16).Value = 1350

# Row 19
$ws.Cells.Item(19, 4).Value = 44460
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 15000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 15500
$ws.Cells.Item(19, 16).Value = 1550

# Row 20
$ws.Cells.Item(20, 4).Value = 44580
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 18000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 19000
$ws.Cells.Item(20, 16).Value = 1900

# Row 21
$ws.Cells.Item(21, 4).Value = 44980
$ws.Cells.Item(21, 10).Value = 400
$ws.Cells.Item(21, 11).Value = 15000
$ws.Cells.Item(21, 12).Value = 16000
$ws.Cells.Item(21, 13).Value = 15500
$ws.Cells.Item(21, 16).Value = 1550

# Row 22
$ws.Cells.Item(22, 4).Value = 44644
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 20000
$ws.Cells.Item(22, 12).Value = 21000
$ws.Cells.Item(22, 13).Value = 20500
$ws.Cells.Item(22, 16).Value = 2050

# Row 23
$ws.Cells.Item(23, 4).Value = 44524
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 20000
$ws.Cells.Item(23, 12).Value = 21000
$ws.Cells.Item(23, 13).Value = 20500
$ws.Cells.Item(23, 16).Value = 2050

# Row 24
$ws.Cells.Item(24, 4).Value = 44904
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 14000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 14500
$ws.Cells.Item(24, 16).Value = 1450

# Row 25
$ws.Cells.Item(25, 4).Value = 44263
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 15000
$ws.Cells.Item(25, 12).Value = 16000
$ws.Cells.Item(25, 13).Value = 15500
$ws.Cells.Item(25, 16).Value = 1550

# Row 26
$ws.Cells.Item(26, 4).Value = 44265
$ws.Cells.Item(26, 10).Value = 200
$ws.Cells.Item(26, 11).Value = 15000
$ws.Cells.Item(26, 12).Value = 16000
$ws.Cells.Item(26, 13).Value = 15500
$ws.Cells.Item(26, 16).Value = 1550

# Row 27
$ws.Cells.Item(27, 4).Value = 44882
$ws.Cells.Item(27, 10).Value = 400
$ws.Cells.Item(27, 11).Value = 15000
$ws.Cells.Item(27, 12).Value = 16000
$ws.Cells.Item(27, 13).Value = 15550
$ws.Cells.Item(27, 16).Value = 1555

# Row 28
$ws.Cells.Item(28, 4).Value = 44988
$ws.Cells.Item(28, 10).Value = 700
$ws.Cells.Item(28, 11).Value = 16000
$ws.Cells.Item(28, 12).Value = 17000
$ws.Cells.Item(28, 13).Value = 16357
$ws.Cells.Item(28, 16).Value = 1636

# Row 29
$ws.Cells.Item(29, 4).Value = 44679
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 19000
$ws.Cells.Item(29, 12).Value = 20000
$ws.Cells.Item(29, 13).Value = 19500
$ws.Cells.Item(29, 16).Value = 1950

# Row 30
$ws.Cells.Item(30, 4).Value = 45008
$ws.Cells.Item(30, 10).Value = 1750
$ws.Cells.Item(30, 11).Value = 18000
$ws.Cells.Item(30, 12).Value = 19000
$ws.Cells.Item(30, 13).Value = 18500
$ws.Cells.Item(30, 16).Value = 1850

# Row 31
$ws.Cells.Item(31, 4).Value = 44694
$ws.Cells.Item(31, 10).Value = 400
$ws.Cells.Item(31, 11).Value = 16000
$ws.Cells.Item(31, 12).Value = 17000
$ws.Cells.Item(31, 13).Value = 16500
$ws.Cells.Item(31, 16).Value = 1650

# Row 32
$ws.Cells.Item(32, 4).Value = 44218
$ws.Cells.Item(32, 10).Value = 320
$ws.Cells.Item(32, 11).Value = 10000
$ws.Cells.Item(32, 12).Value = 11000
$ws.Cells.Item(32, 13).Value = 10500
$ws.Cells.Item(32, 16).Value = 1050

# Row 33
$ws.Cells.Item(33, 4).Value = 44428
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 15000
$ws.Cells.Item(33, 12).Value = 16000
$ws.Cells.Item(33, 13).Value = 15500
$ws.Cells.Item(33, 16).Value = 1550

# Row 34
$ws.Cells.Item(34, 4).Value = 44714
$ws.Cells.Item(34, 10).Value = 400
$ws.Cells.Item(34, 11).Value = 19000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 13).Value = 19500
$ws.Cells.Item(34, 16).Value = 1950

# Row 35
$ws.Cells.Item(35, 4).Value = 44925
$ws.Cells.Item(35, 10).Value = 250
$ws.Cells.Item(35, 11).Value = 14000
$ws.Cells.Item(35, 12).Value = 15000
$ws.Cells.Item(35, 13).Value = 14600
$ws.Cells.Item(35, 16).Value = 1460

# Row 36
$ws.Cells.Item(36, 4).Value = 44847
$ws.Cells.Item(36, 10).Value = 400
$ws.Cells.Item(36, 11).Value = 16000
$ws.Cells.Item(36, 12).Value = 17000
$ws.Cells.Item(36, 13).Value = 16500
$ws.Cells.Item(36, 16).Value = 1650

# Row 37
$ws.Cells.Item(37, 4).Value = 44777
$ws.Cells.Item(37, 10).Value = 200
$ws.Cells.Item(37, 11).Value = 24000
$ws.Cells.Item(37, 12).Value = 25000
$ws.Cells.Item(37, 13).Value = 24500
$ws.Cells.Item(37, 16).Value = 2450

# Row 38
$ws.Cells.Item(38, 4).Value = 44441
$ws.Cells.Item(38, 10).Value = 300
$ws.Cells.Item(38, 11).Value = 15000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 15500
$ws.Cells.Item(38, 16).Value = 1550

# Row 39
$ws.Cells.Item(39, 4).Value = 44727
$ws.Cells.Item(39, 10).Value = 400
$ws.Cells.Item(39, 11).Value = 18000
$ws.Cells.Item(39, 12).Value = 19000
$ws.Cells.Item(39, 13).Value = 18500
$ws.Cells.Item(39, 16).Value = 1850

# Row 40
$ws.Cells.Item(40, 4).Value = 44918
$ws.Cells.Item(40, 10).Value = 200
$ws.Cells.Item(40, 11).Value = 12000
$ws.Cells.Item(40, 12).Value = 13000
$ws.Cells.Item(40, 13).Value = 12250
$ws.Cells.Item(40, 16).Value = 1225
